$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.286676049232483
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 2.196618556976318
$ws.Range("D1").Value = 1.291357278823853
$ws.Range("E1").Value = 0.9157800674438477
